$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Thresholds in Database")

$newVersionString = "IndicatorQuantiles.R, Git Commit ID: db49f0f869e1f5a8558dc746458075a467cf2c41"
$newSourceString = "Database_Thresholds_details.xlsx, Git Commit ID: 17b6a0f858dccbb28fc8ab3fe179e7fa731e5996"

# Column W (ScriptLatestRunVersion) is updated on every data row (4-92).
for ($r = 4; $r -le 92; $r++) {
    $ws.Cells.Item($r, 23).Value2 = $newVersionString   # column W
}

# Column U (QuantileSource) is updated only on these specific rows.
$sourceRows = @(29, 41, 62, 63, 78)
foreach ($r in $sourceRows) {
    $ws.Cells.Item($r, 21).Value2 = $newSourceString    # column U
}
